$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the 2020 column (N) mirroring the style/format of the existing 2019 column (M)
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2020

$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 2.1

# Reflect the selection that was active when the workbook was last saved
$ws.Activate()
$ws.Range("N9").Select()
